# Update the "dSF" column (F) values for several rows in the sheet.
# These values were repulled/recalculated from source data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = 2
    3  = 9
    6  = -1
    7  = -3
    8  = -2
    11 = -10
    12 = 7
    14 = -2
    15 = -3
    16 = 0
    17 = 2
    18 = -1
    19 = 1
    20 = -4
    21 = -1
    23 = -3
    24 = -1
    25 = -1
    27 = 2
    28 = 3
    29 = 1
    30 = -2
    31 = -1
    33 = -3
    34 = 1
    35 = -4
    36 = 3
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
